$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 10
$ws.Range("E12").Value = "error_message"
$ws.Range("F12").Value = "現在、あなたの質問に関するデータはありません。具体的な質問をしていただければ、よりお手伝いできるかと思います。"
$ws.Range("G12").Value = "ユーザーの質問に対する該当データが存在しない場合に表示されるエラーメッセージ。より具体的な質問を促すために使用。"
